# consultaPrecos/produtos.xlsx - "Abre site e pesquisa produto"
#
# The price column header is renamed and a new blank spacer row is
# inserted under the title row (row 2), pushing the product rows down
# by one. The selection cursor is left on B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2 (shifts "Monitor..."/"Mouse..." rows down).
$ws.Rows.Item(2).Insert()

# Rename the B1 header from "Precos" to "Preços/data".
$ws.Range("B1").Value = "Preços/data"

# Leave the selection on B8, matching the saved cursor position.
$ws.Range("B8").Select() | Out-Null
